$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-6 hold a "mahmoud / 1 / 29 / 16-9-2025" record.
# Append the same record in rows 7-16 (matching style of row 2).
for ($r = 7; $r -le 16; $r++) {
    $ws.Cells.Item($r, 1).Value = "mahmoud"
    $ws.Cells.Item($r, 2).Value = 1
    $ws.Cells.Item($r, 3).Value = 29
    $ws.Cells.Item($r, 4).Value = "16-9-2025"

    # Match the formatting/style used by the existing data rows
    $ws.Range("A2:D2").Copy() | Out-Null
    $ws.Range("A" + $r + ":D" + $r).PasteSpecial(-4122) | Out-Null
}

$excel.CutCopyMode = 0

$ws.Range("F11").Select() | Out-Null
